# Apply the edits described by the commit diff:
#  - Rename "Sheet1" to "Demo"
#  - Add two new rows of data: A3 = "Active Sync2", A4 = "Cat"
#  - Leave the active selection on A3 (as captured in the saved file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Demo"

$ws.Range("A4").Value = "Cat"
$ws.Range("A3").Value = "Active Sync2"

$ws.Range("A3").Select()
